# Generate Report for Archive
# - Flip "Ready for handoff" status cells over to "In Translation"
#   (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3 all share the same status string)
# - Shrink the now-narrower "Status" columns to match the new content

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn / de-de status columns (E, F) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = $newStatus
$overview.Range("E1:F1").EntireColumn.ColumnWidth = 12.5

# --- zh-cn sheet: Status column (C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C3").Value = $newStatus
$zhcn.Range("C1").EntireColumn.ColumnWidth = 12.5

# --- de-de sheet: Status column (C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C3").Value = $newStatus
$dede.Range("C1").EntireColumn.ColumnWidth = 12.5
